$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.006.41'
$ws.Range("D3").Value = '1.740.29'
$ws.Range("E3").Value = '  +5.16%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").Value = '''228.78'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.41%  '
$ws.Range("D6").Value = '''0.5447'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.88%  '
$ws.Range("D7").Value = '''1.003'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.16%  '
$ws.Range("D8").Value = '''0.2772'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +4.06%  '
$ws.Range("D9").Value = '''0.06758'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +6.29%  '
$ws.Range("D10").Value = '''21.71'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.89%  '
$ws.Range("D11").Value = '''0.07794'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.91%  '
$ws.Range("D12").Value = '''4.705'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.35%  '
$ws.Range("D13").Value = '1.736.74'
$ws.Range("E13").Value = '  +3.78%  '
$ws.Range("D14").Value = '1.981.41'
$ws.Range("E14").Value = '  +5.22%  '
$ws.Range("D15").Value = '''0.5998'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +6.62%  '
$ws.Range("D16").Value = '0.0₅8414'
$ws.Range("E16").Value = '  +1.98%  '
$ws.Range("D17").Value = '''69.27'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +5.95%  '
$ws.Range("D18").Value = '27.979.09'
$ws.Range("E18").Value = '  +6.70%  '
$ws.Range("D19").Value = '''225.12'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +17.25%  '
$ws.Range("D20").Value = '''4.848'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.25%  '
$ws.Range("E21").Value = '  -0.14%  '
$ws.Range("D22").Value = '''10.93'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +5.08%  '
$ws.Range("D23").Value = '''6.233'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.81%  '
$ws.Range("D24").Value = '''1.004'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.14%  '
$ws.Range("D25").Value = '''146.35'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.86%  '
$ws.Range("D26").Value = '''0.1252'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.92%  '
$ws.Range("D27").Value = '''7.470'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.74%  '
$ws.Range("D28").Value = '''17.09'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +7.18%  '
$ws.Range("D29").Value = '''1.648'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +8.79%  '
$ws.Range("D30").Value = '''0.05690'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.09%  '
$ws.Range("D31").Value = '''1.319'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.30%  '
$ws.Range("D32").Value = '''3.716'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +6.04%  '
$ws.Range("D33").Value = '''3.534'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.20%  '
$ws.Range("D34").Value = '''1.673'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +5.56%  '
$ws.Range("D35").Value = '''0.9830'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.04%  '
$ws.Range("D36").Value = '''2.860'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.12%  '
$ws.Range("D37").Value = '''2.452'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.59%  '
$ws.Range("D38").Value = '''0.5961'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.52%  '
$ws.Range("D39").Value = '''0.01677'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.81%  '
$ws.Range("D40").Value = '''6.004'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.06%  '
$ws.Range("D41").Value = '''0.8494'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.87%  '
$ws.Range("D42").Value = '1.048.43'
$ws.Range("E42").Value = '  +4.00%  '
$ws.Range("D43").Value = '''1.003'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.08%  '
$ws.Range("D44").Value = '''102.24'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.28%  '
$ws.Range("D45").Value = '1.885.47'
$ws.Range("E45").Value = '  +5.06%  '
$ws.Range("D46").Value = '0.0₈117'
$ws.Range("E46").Value = '  +14.83%  '
$ws.Range("D47").Value = '''60.08'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.88%  '
$ws.Range("D48").Value = '''8.312'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.01%  '
$ws.Range("D49").Value = '''0.4431'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.90%  '
$ws.Range("D50").Value = '''1.005'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.23%  '
$ws.Range("D51").Value = '''0.05320'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.44%  '
